$wb = $excel.ActiveWorkbook
$wsOpt = $wb.Worksheets.Item("optimization_parameters")

# --- optimization_parameters sheet restructuring ---
# Row 1: drop the extra "value" labels in C1:F1 (only A1/B1 remain)
$wsOpt.Range("C1:F1").ClearContents()

# Insert a new row at 9 (shifts old rows 9-17 down to 10-18)
$wsOpt.Rows.Item(9).Insert()

# The old "Deletion" row (now at row 17 after the insert) is removed entirely,
# shifting the final "simulation_timepoints" row back up to row 17
$wsOpt.Rows.Item(17).Delete()

# Rename "Model" (A8) to "production_function"
$wsOpt.Range("A8").Value = "production_function"

# Populate the newly inserted row 9 with the "L_curve" parameter
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 0
$wsOpt.Range("B9").NumberFormat = "0.00E+00"

# --- Active sheet / selection changes ---
# Make "optimization_parameters" the active tab with C1:F1 selected
$wsOpt.Activate()
$wsOpt.Range("C1:F1").Select()
